# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet "NewLoanInput": rows 15 & 18 carried a stale explicit row height from
# a previously longer wrapped label; re-autofit them back to the sheet's
# normal height, then leave the selection where the user last clicked.
# ---------------------------------------------------------------------------
$wsLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsLoanInput.Activate()
$wsLoanInput.Rows.Item(15).EntireRow.AutoFit()
$wsLoanInput.Rows.Item(18).EntireRow.AutoFit()
$wsLoanInput.Range("A11").Select()

# ---------------------------------------------------------------------------
# Sheet "Summary": just a cursor/selection move.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule": the loan now runs a full 12 instalments instead
# of 6, so six more repayment rows are appended (9-14), and every data row
# gains an extra "0" column (O) that was missing before.
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Activate()

# Rows 3-8 each get a new "0" in column O (format copied from column N so it
# keeps the same style as the rest of the row).
$wsSchedule.Range("N3:N8").Copy()
$wsSchedule.Range("O3:O8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$wsSchedule.Range("O3").Value = 0
$wsSchedule.Range("O4").Value = 0
$wsSchedule.Range("O5").Value = 0
$wsSchedule.Range("O6").Value = 0
$wsSchedule.Range("O7").Value = 0
$wsSchedule.Range("O8").Value = 0

# Row 2 (the disbursement row) only needs a blank trailing cell added.
$wsSchedule.Range("O2").Copy()
$wsSchedule.Range("P2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Use row 8 as the template (formats) for the six new instalment rows.
$wsSchedule.Range("A8:P8").Copy()
$wsSchedule.Range("A9:P14").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$scheduleRows = @(
    @{ Row=9;  A=7;  B=31; C=42217; F=834.99; G=4338.57; H=52.73; K=887.72; P=887.72 },
    @{ Row=10; A=8;  B=31; C=42248; F=843.5;  G=3495.07; H=44.22; K=887.72; P=887.72 },
    @{ Row=11; A=9;  B=30; C=42278; F=853.25; G=2641.82; H=34.47; K=887.72; P=887.72 },
    @{ Row=12; A=10; B=31; C=42309; F=860.8;  G=1781.02; H=26.92; K=887.72; P=887.72 },
    @{ Row=13; A=11; B=30; C=42339; F=870.15; G=910.87;  H=17.57; K=887.72; P=887.72 },
    @{ Row=14; A=12; B=31; C=42370; F=910.87; G=0;       H=9.28;  K=920.15; P=920.15 }
)

foreach ($r in $scheduleRows) {
    $row = $r.Row
    $wsSchedule.Range("A$row").Value = $r.A
    $wsSchedule.Range("B$row").Value = $r.B
    $wsSchedule.Range("C$row").Value = $r.C
    $wsSchedule.Range("F$row").Value = $r.F
    $wsSchedule.Range("G$row").Value = $r.G
    $wsSchedule.Range("H$row").Value = $r.H
    $wsSchedule.Range("I$row").Value = 0
    $wsSchedule.Range("J$row").Value = 0
    $wsSchedule.Range("K$row").Value = $r.K
    $wsSchedule.Range("L$row").Value = 0
    $wsSchedule.Range("M$row").Value = 0
    $wsSchedule.Range("N$row").Value = 0
    $wsSchedule.Range("O$row").Value = 0
    $wsSchedule.Range("P$row").Value = $r.P
}

# The last two instalments' running-balance column (G) renders in the plain
# style used elsewhere on the row instead of the "#,##0.00" style used
# further up the schedule - copy H13/H14's format (style 12) onto G13/G14.
$wsSchedule.Range("H13").Copy()
$wsSchedule.Range("G13").PasteSpecial($xlPasteFormats)
$wsSchedule.Range("H14").Copy()
$wsSchedule.Range("G14").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$wsSchedule.Range("G13").Value = 910.87
$wsSchedule.Range("G14").Value = 0

$wsSchedule.Rows.Item(15).Select()

# ---------------------------------------------------------------------------
# Sheet "Transactions": the running transaction-id numbering was corrected.
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
$wsTransactions.Range("A2").Value = 87
$wsTransactions.Range("A3").Value = 86
$wsTransactions.Range("A4").Value = 85
$wsTransactions.Range("A5").Value = 84
$wsTransactions.Range("A2:L5").Select()
